# Apply the Tue Oct 22 04:47:51 UTC 2024 cryptos list refresh (GitHub Actions bot edit).
# Updates the Price (column D) and Volume(1h) (column E) columns for rows 2-51
# of the active worksheet to their freshly-scraped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some refreshed Price cells look like plain numbers (e.g. '598.14', '19.40').
# The source data stores Price/Volume as literal text, so force those specific
# cells to Text format before writing -- otherwise Excel's COM layer would
# auto-coerce them to numbers (and e.g. drop the trailing zero in '19.40').
$textForceCells = @("D5", "D6", "D8", "D12", "D14", "D19", "D21", "D23", "D24", "D27", "D29", "D31", "D32", "D34", "D39", "D40", "D42", "D46", "D49", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '67.427.31'
$ws.Range("E2").Value = '  -2.15%  '
$ws.Range("D3").Value = '2.641.33'
$ws.Range("E3").Value = '  -3.31%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '598.14'
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("D6").Value = '167.04'
$ws.Range("E6").Value = '  -1.39%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '0.543'
$ws.Range("E8").Value = '  -0.68%  '
$ws.Range("D9").Value = '2.640.89'
$ws.Range("E9").Value = '  -3.28%  '
$ws.Range("E10").Value = '  -0.91%  '
$ws.Range("E11").Value = '  +1.40%  '
$ws.Range("D12").Value = '0.364'
$ws.Range("E12").Value = '  -1.45%  '
$ws.Range("E13").Value = '  -2.14%  '
$ws.Range("D14").Value = '27.95'
$ws.Range("E14").Value = '  -2.68%  '
$ws.Range("D15").Value = '3.122.14'
$ws.Range("E15").Value = '  -3.31%  '
$ws.Range("E16").Value = '  -3.49%  '
$ws.Range("D17").Value = '67.317.47'
$ws.Range("E17").Value = '  -2.09%  '
$ws.Range("D18").Value = '2.643.95'
$ws.Range("E18").Value = '  -3.16%  '
$ws.Range("D19").Value = '11.87'
$ws.Range("E19").Value = '  -0.60%  '
$ws.Range("E20").Value = '  +2.19%  '
$ws.Range("D21").Value = '362.95'
$ws.Range("E21").Value = '  -2.90%  '
$ws.Range("E22").Value = '  -3.11%  '
$ws.Range("D23").Value = '4.78'
$ws.Range("E23").Value = '  -3.92%  '
$ws.Range("D24").Value = '10.91'
$ws.Range("E24").Value = '  +8.08%  '
$ws.Range("E25").Value = '  -5.85%  '
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("D27").Value = '70.87'
$ws.Range("E27").Value = '  -3.78%  '
$ws.Range("D28").Value = '2.779.94'
$ws.Range("E28").Value = '  -3.23%  '
$ws.Range("D29").Value = '0.0000102'
$ws.Range("E29").Value = '  -3.50%  '
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("D31").Value = '553.99'
$ws.Range("E31").Value = '  -5.52%  '
$ws.Range("D32").Value = '8.06'
$ws.Range("E32").Value = '  -3.23%  '
$ws.Range("E33").Value = '  -4.03%  '
$ws.Range("D34").Value = '1.92'
$ws.Range("E34").Value = '  -1.95%  '
$ws.Range("E35").Value = '  +0.31%  '
$ws.Range("E37").Value = '  -5.37%  '
$ws.Range("E38").Value = '  -2.73%  '
$ws.Range("D39").Value = '19.40'
$ws.Range("E39").Value = '  -3.20%  '
$ws.Range("D40").Value = '0.372'
$ws.Range("E40").Value = '  -2.57%  '
$ws.Range("E41").Value = '  -5.32%  '
$ws.Range("D42").Value = '5.26'
$ws.Range("E42").Value = '  -4.38%  '
$ws.Range("E43").Value = '  -0.37%  '
$ws.Range("E44").Value = '  -5.23%  '
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").Value = '40.15'
$ws.Range("E46").Value = '  -2.24%  '
$ws.Range("D47").Value = '0.0₆0300'
$ws.Range("E47").Value = '  -3.91%  '
$ws.Range("E48").Value = '  -1.72%  '
$ws.Range("D49").Value = '153.55'
$ws.Range("E49").Value = '  -1.75%  '
$ws.Range("D50").Value = '3.88'
$ws.Range("E50").Value = '  -2.56%  '
$ws.Range("E51").Value = '  -4.08%  '
